# "Add 3 more tests for Login page" - adds 3 new data rows (IDs 2-4) to the
# DataSet sheet of the login-tests workbook, re-styles the existing/new
# "Email" cells as hyperlinks (mailto: links), and widens the Key/Email
# columns to fit the new, longer test-case names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string values -------------------------------------------
# These four brand-new strings must be *first written* in this exact order
# (LoginWithoutPassword, LoginWithoutEmailAndPassword,
#  LoginWithExistingUserButWrongPassword, Wrong) so that the workbook's
# shared-string table indices land the same way they do in the target file
# (existing_mail@sm.bg is already shared from row 2 and does not need to be
# re-added).
$ws.Range("B4").Value = "LoginWithoutPassword"
$ws.Range("B3").Value = "LoginWithoutEmailAndPassword"
$ws.Range("B5").Value = "LoginWithExistingUserButWrongPassword"
$ws.Range("D5").Value = "Wrong"

# --- Row 3: LoginWithoutEmailAndPassword (Key only) ----------------------
$ws.Range("A3").Value = 2

# --- Row 4: LoginWithoutPassword (Key + Email) ----------------------------
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = "existing_mail@sm.bg"

# --- Row 5: LoginWithExistingUserButWrongPassword (Key + Email + Password)
$ws.Range("A5").Value = 4
$ws.Range("C5").Value = "existing_mail@sm.bg"

# --- Hyperlink-ify every Email cell (existing row 2 + the two new ones) --
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:existing_mail@sm.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:existing_mail@sm.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:existing_mail@sm.bg") | Out-Null

# --- Column widths: Key/Email columns now "best fit" the longer content --
$ws.Range("B:B").ColumnWidth = 35.498697916666664
$ws.Range("C:C").ColumnWidth = 17.498697916666668

# --- Leave the selection where the author's saved file shows it ----------
$ws.Range("B9").Select() | Out-Null
